$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.250.88'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.591.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.245'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0853'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.817.02'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.590.69'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.83'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.248.88'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0725'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.28'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('E24').Value = '  -3.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.10'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.412.02'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.23%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('E36').Value = '  -1.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.577'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.968'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.94%  '
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.728.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.93'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0501'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0951'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.54%  '
$ws.Range('E51').Value = '  +0.10%  '
